# Weekly update: insert a new week's worth of data (3 rows: Primera / Segunda /
# Tercera) at the top of the existing data block for this product/market, pushing
# all subsequent rows down by three rows. This mirrors how each workbook in this
# series prepends the latest week's prices above the historical records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows right before the first data row of the block (row 485),
# shifting the existing data (previously rows 485:605) down to rows 488:608.
$ws.Range("A485:R487").EntireRow.Insert()

# Values shared by every row of this market/product (unchanged from the rest of
# the sheet).
$mercado = "Agrícola del Norte S.A. de Arica"
$region = "Arica y Parinacota"
$categoriaId = 100112020
$categoria = "Tomate"
$variedad = "Larga vida"
$unidad = "`$/caja 10 kilos"
$origen = "Región de Arica y Parinacota"
$clasificacion = "Hortaliza"

# New week's data (fecha = 44551) for the three quality grades.
$newData = @(
    @{ Row = 485; Calidad = "Primera"; Volumen = 250; PMin = 3000; PMax = 3500; PProm = 3250; PKg = 325 },
    @{ Row = 486; Calidad = "Segunda"; Volumen = 270; PMin = 2500; PMax = 3000; PProm = 2750; PKg = 275 },
    @{ Row = 487; Calidad = "Tercera"; Volumen = 300; PMin = 2000; PMax = 2500; PProm = 2250; PKg = 225 }
)

foreach ($entry in $newData) {
    $r = $entry.Row

    $ws.Cells.Item($r, 1).Value = 1
    $ws.Cells.Item($r, 2).Value = $mercado
    $ws.Cells.Item($r, 3).Value = $region
    $ws.Cells.Item($r, 4).Value = 44551
    $ws.Cells.Item($r, 5).Value = 15
    $ws.Cells.Item($r, 6).Value = $categoriaId
    $ws.Cells.Item($r, 7).Value = $categoria
    $ws.Cells.Item($r, 8).Value = $variedad
    $ws.Cells.Item($r, 9).Value = $entry.Calidad
    $ws.Cells.Item($r, 10).Value = $entry.Volumen
    $ws.Cells.Item($r, 11).Value = $entry.PMin
    $ws.Cells.Item($r, 12).Value = $entry.PMax
    $ws.Cells.Item($r, 13).Value = $entry.PProm
    $ws.Cells.Item($r, 14).Value = $unidad
    $ws.Cells.Item($r, 15).Value = $origen
    $ws.Cells.Item($r, 16).Value = $entry.PKg
    $ws.Cells.Item($r, 17).Value = 10
    $ws.Cells.Item($r, 18).Value = $clasificacion
}
